$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 31: 2026-06-21, Lancamento Glow Facial / Ollie CL / Lançamento / G / Em Dev / Previsão ---
$ws.Range("A31").Value = 46194
$ws.Range("A31").NumberFormat = 'yyyy\-mm\-dd\ hh:mm:ss'

$ws.Range("B31").Value = "Lancamento Glow Facial"
$ws.Range("B31").Font.Color = $ws.Range("B30").Font.Color

$ws.Range("C31").Value = "Ollie CL"
$ws.Range("C31").Font.Color = $ws.Range("B30").Font.Color

$ws.Range("D31").Value = "Lançamento"
$ws.Range("D31").Font.Color = $ws.Range("B30").Font.Color

$ws.Range("E31").Value = "G"
$ws.Range("E31").Font.Color = $ws.Range("B30").Font.Color

$ws.Range("F31").Value = "🟡 Em Dev"
$ws.Range("F31").Font.Color = $ws.Range("B30").Font.Color

$ws.Range("G31").Value = "Previsão"
$ws.Range("G31").Font.Color = $ws.Range("B30").Font.Color

# --- Row 32: 2026-06-22, Lancamento Glow Facial / Ollie CL / Lançamento / G / Em Dev / Previsão ---
$ws.Range("A32").Value = 46195
$ws.Range("A32").NumberFormat = 'yyyy\-mm\-dd\ hh:mm:ss'

$ws.Range("B32").Value = "Lancamento Glow Facial"
$ws.Range("B32").Font.Color = $ws.Range("B30").Font.Color

$ws.Range("C32").Value = "Ollie CL"
$ws.Range("C32").Font.Color = $ws.Range("B30").Font.Color

$ws.Range("D32").Value = "Lançamento"
$ws.Range("D32").Font.Color = $ws.Range("B30").Font.Color

$ws.Range("E32").Value = "G"
$ws.Range("E32").Font.Color = $ws.Range("B30").Font.Color

$ws.Range("F32").Value = "🟡 Em Dev"
$ws.Range("F32").Font.Color = $ws.Range("B30").Font.Color

$ws.Range("G32").Value = "Previsão"
$ws.Range("G32").Font.Color = $ws.Range("B30").Font.Color

# --- Update selection / active cell to match the saved view state ---
$excel.Goto($ws.Range("A31:G32"))
